# Pruebas y modificaciones de servicios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string used in C16: "tipoCategorias" -> "tiposCategorias"
$ws.Range("C16").Value = "tiposCategorias"

# New numeric data added to columns P, Q, R (service test/measurement figures)
$ws.Range("Q18").Value = 5
$ws.Range("R18").Value = 5

$ws.Range("P19").Value = 5
$ws.Range("Q19").Value = 8
$ws.Range("R19").Value = 2

$ws.Range("P20").Value = 10
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = -3

$ws.Range("R21").Value = 7

$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5

$ws.Range("P25").Value = 10
$ws.Range("Q25").Value = 8

$ws.Range("P26").Value = 5
$ws.Range("Q26").Value = 3

$ws.Range("Q30").Value = 5

$ws.Range("P31").Value = 10
$ws.Range("Q31").Value = 4

$ws.Range("P32").Value = 5

$ws.Range("Q34").Value = 5

$ws.Range("P35").Value = 10
$ws.Range("Q35").Value = 9

$ws.Range("P36").Value = 15

# Update the active selection to match the new working cell
$ws.Range("G27").Select()
